$d = $word.ActiveDocument

# Locate the paragraph containing "27/10-2025" (the heading right before
# where the new diary entry needs to go).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*27/10-2025*") {
        $target = $p
        break
    }
}

$insPoint = $d.Range($target.Range.End, $target.Range.End)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>I dag har jeg lavet alle mappers til mine DTO' + [char]8217 + 'er og Modeller, samt lavet min passwordhaser s' + [char]229 + ' vi kan sikre brugernes informationer bedre.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insPoint.InsertXML($xmlFrag)

# InsertXML leaves a stray empty paragraph behind (an artifact of how the
# trailing fragment paragraph merges with the following "28/10-2025"
# paragraph) -- remove it so the document matches the intended structure.
$stray = $target.Next().Next()
if ($stray.Range.Text -eq "") {
    $stray.Range.Delete()
}

Write-Host "done"
